{"js": "// Prepend \"Design: \" to the start of each feedback answer in the\n// Q&A table (the \"List Bullet\" styled paragraph inside the answer\n// column of every data row, skipping the header row).\nconst table = context.document.body.tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst prefix = \"Design: \";\n\nfor (let r = 1; r < table.rowCount; r++) {\n  const cell = table.getCell(r, 1);\n  const paragraphs = cell.body.paragraphs;\n  paragraphs.load(\"items/text,items/style\");\n  await context.sync();\n\n  for (const paragraph of paragraphs.items) {\n    if (paragraph.style === \"List Bullet\" && !paragraph.text.startsWith(prefix)) {\n      paragraph.insertText(prefix, \"Start\");\n    }\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n$prefix = \"Design: \"\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Style.NameLocal -eq \"List Bullet\") {\n        $r = $p.Range\n        if ($r.Text.IndexOf($prefix) -ne 0) {\n            $r.InsertBefore($prefix)\n        }\n    }\n}\n"}
